# "Generate Report for Handoff"
# The c6c61c5a\...md file (row 3 in each sheet) has moved from
# "Handed back: in sync with en-US" to "Ready for handoff", with a new
# handoff timestamp and an error message noting the handback file is stale.

$wb = $excel.ActiveWorkbook

$statusReady      = "Ready for handoff"
$handoffDateTime  = "2016-08-28 02:48:10"
$handoffDateTimeZh = "2016-08-28 02:48:06"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/566a9f252b58d8bbab1951b3ec466c921a550977/e2e/c6c61c5a-1579-4e80-85c7-39e8865daba8.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f468c6462967d0b76061cc09bb622781848cc684/e2e/c6c61c5a-1579-4e80-85c7-39e8865daba8.md."

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusReady
$wsOverview.Range("F3").Value = $statusReady
$wsOverview.Range("G3").Value = $handoffDateTime

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = $statusReady
$wsZh.Range("H3").Value = $handoffDateTimeZh
$wsZh.Range("P3").Value = $errorDetail
$wsZh.Columns.Item(16).ColumnWidth = 39.17

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = $statusReady
$wsDe.Range("H3").Value = $handoffDateTime
$wsDe.Range("P3").Value = $errorDetail
$wsDe.Columns.Item(16).ColumnWidth = 39.17
